# Rubric change and tweaks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Milestone Completed" (E) and "Complete(X)" (F) markers ---

# Milestone III only entries
$ws.Range("E6").Value = "III"
$ws.Range("E20").Value = "III"
$ws.Range("E29").Value = "III"
$ws.Range("E30").Value = "III"
$ws.Range("E31").Value = "III"
$ws.Range("E53").Value = "III"
$ws.Range("F53").Value = "X"
$ws.Range("E54").Value = "III"
$ws.Range("E74").Value = "III"

# Milestone II, marked complete with X
$ws.Range("E25").Value = "II"
$ws.Range("F25").Value = "X"
$ws.Range("E26").Value = "II"
$ws.Range("F26").Value = "X"
$ws.Range("E27").Value = "II"
$ws.Range("F27").Value = "X"
$ws.Range("E28").Value = "II"
$ws.Range("F28").Value = "X"
$ws.Range("E35").Value = "II"
$ws.Range("F35").Value = "X"
$ws.Range("E39").Value = "II"
$ws.Range("F39").Value = "X"

# Already "II" entries, now marked complete with X
$ws.Range("F62").Value = "X"
$ws.Range("F64").Value = "X"

# Carry-over completion flags
$ws.Range("D83").Value = "X"
$ws.Range("E83").Value = "X"
$ws.Range("D84").Value = "X"
$ws.Range("E84").Value = "X"

# --- Update selection / view state ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("F74").Select()
